$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two species-observation rows (row 2 <-> row 3) for the
# columns whose values actually differ between the rows: A (Id), B
# (Taxonsorteringsordning), D (Rodlistade), E (TaxonId), F (Artnamn),
# G (Vetenskapligt namn), H (Auktor) and AH (Biotop).

$cols = @("A","B","D","E","F","G","H","AH")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}
